$wb = $excel.ActiveWorkbook

# --- Sheet 1: Significant Components ---
$ws1 = $wb.Worksheets.Item("Significant Components")
$ws1.Range("C2").Value2 = "['QSERV' 'QFHH' 'QHISPC' 'QEDLESHI' 'QNOHLTH' 'QESL' 'PPUNIT' 'QEXTRCT'`n 'PERCAP']"
$ws1.Range("C3").Value2 = "['QRICH' 'PERCAP' 'MDHSEVAL']"
$ws1.Range("C4").Value2 = "['PPUNIT' 'QNOAUTO' 'QPOVTY' 'QRENTER' 'QFAM']"
$ws1.Range("C5").Value2 = "['QAGEDEP' 'QFEMLBR' 'QFEMALE']"
$ws1.Range("C6").Value2 = "['QAGEDEP' 'QSSBEN' 'MEDAGE']"

# --- Sheet 2: Loading Factors ---
$ws2 = $wb.Worksheets.Item("Loading Factors")
$ws2.Range("A2").Value2 = "QSERV"
$ws2.Range("B2").Value2 = 0.5381204751319512
$ws2.Range("C2").Value2 = 0.3662181495109431
$ws2.Range("D2").Value2 = 0.2976695030656286
$ws2.Range("E2").Value2 = 0.001243119697268828
$ws2.Range("F2").Value2 = -0.169403188940627
$ws2.Range("A3").Value2 = "QFHH"
$ws2.Range("B3").Value2 = 0.6130097831512937
$ws2.Range("C3").Value2 = 0.2399781591693043
$ws2.Range("D3").Value2 = 0.1790883841059889
$ws2.Range("E3").Value2 = 0.262013816775974
$ws2.Range("F3").Value2 = -0.09687619950102999
$ws2.Range("A4").Value2 = "QHISPC"
$ws2.Range("B4").Value2 = 0.8388032637219848
$ws2.Range("C4").Value2 = 0.3352021753304731
$ws2.Range("D4").Value2 = 0.1601474956817818
$ws2.Range("E4").Value2 = -0.06922549985931019
$ws2.Range("F4").Value2 = -0.09574193465063523
$ws2.Range("A5").Value2 = "QEDLESHI"
$ws2.Range("B5").Value2 = 0.8855973614226529
$ws2.Range("C5").Value2 = 0.2225893013927557
$ws2.Range("D5").Value2 = 0.2488821997052932
$ws2.Range("E5").Value2 = -0.1306950246605276
$ws2.Range("F5").Value2 = 0.01576000589942111
$ws2.Range("A6").Value2 = "QNOHLTH"
$ws2.Range("B6").Value2 = 0.6164840041304507
$ws2.Range("C6").Value2 = 0.4453383193800906
$ws2.Range("D6").Value2 = 0.3253808165621783
$ws2.Range("E6").Value2 = -0.2013657286674192
$ws2.Range("F6").Value2 = -0.08210583967820179
$ws2.Range("A7").Value2 = "QESL"
$ws2.Range("B7").Value2 = 0.8069195261612075
$ws2.Range("C7").Value2 = 0.1431421451574245
$ws2.Range("D7").Value2 = 0.2282344082060187
$ws2.Range("E7").Value2 = -0.2351410650504814
$ws2.Range("F7").Value2 = -0.02150689556929739
$ws2.Range("A8").Value2 = "PPUNIT"
$ws2.Range("B8").Value2 = 0.7902275169365798
$ws2.Range("C8").Value2 = -0.02455912049023291
$ws2.Range("D8").Value2 = -0.4600798249827534
$ws2.Range("E8").Value2 = 0.06838368016836407
$ws2.Range("F8").Value2 = -0.1573190495954288
$ws2.Range("A9").Value2 = "QEXTRCT"
$ws2.Range("B9").Value2 = 0.7251699886427621
$ws2.Range("C9").Value2 = 0.1597478160462525
$ws2.Range("D9").Value2 = 0.1377359028071791
$ws2.Range("E9").Value2 = -0.2661556136340222
$ws2.Range("F9").Value2 = 0.05427315761957422
$ws2.Range("A10").Value2 = "QRICH"
$ws2.Range("B10").Value2 = 0.1589560336557637
$ws2.Range("C10").Value2 = 0.8582481198084332
$ws2.Range("D10").Value2 = 0.3041837810343528
$ws2.Range("E10").Value2 = 0.004009562688229646
$ws2.Range("F10").Value2 = -0.1338394787664548
$ws2.Range("A11").Value2 = "PERCAP"
$ws2.Range("B11").Value2 = 0.4725767448997422
$ws2.Range("C11").Value2 = 0.7137769933882232
$ws2.Range("D11").Value2 = 0.23360964660658
$ws2.Range("E11").Value2 = 0.00532654035318215
$ws2.Range("F11").Value2 = -0.2730791446221674
$ws2.Range("A12").Value2 = "MDHSEVAL"
$ws2.Range("B12").Value2 = 0.3726973484607343
$ws2.Range("C12").Value2 = 0.8017743548999904
$ws2.Range("D12").Value2 = 0.08260669463154222
$ws2.Range("E12").Value2 = 0.03067618588106992
$ws2.Range("F12").Value2 = -0.06162804312041255
$ws2.Range("A13").Value2 = "QNOAUTO"
$ws2.Range("B13").Value2 = 0.1911459588519667
$ws2.Range("C13").Value2 = 0.09747169108633319
$ws2.Range("D13").Value2 = 0.6356287686922173
$ws2.Range("E13").Value2 = -0.006055162496571286
$ws2.Range("F13").Value2 = -0.03646340287378295
$ws2.Range("A14").Value2 = "QPOVTY"
$ws2.Range("B14").Value2 = 0.4090324292043737
$ws2.Range("C14").Value2 = 0.1933683399464961
$ws2.Range("D14").Value2 = 0.4940719683132042
$ws2.Range("E14").Value2 = -0.02947174366382025
$ws2.Range("F14").Value2 = -0.3383025141107034
$ws2.Range("A15").Value2 = "QRENTER"
$ws2.Range("B15").Value2 = -0.007034393176870034
$ws2.Range("C15").Value2 = 0.2082000072585756
$ws2.Range("D15").Value2 = 0.7800125845244715
$ws2.Range("E15").Value2 = -0.1101264162299841
$ws2.Range("F15").Value2 = -0.4193040220698699
$ws2.Range("A16").Value2 = "QFAM"
$ws2.Range("B16").Value2 = 0.2493636377250857
$ws2.Range("C16").Value2 = 0.249820823376892
$ws2.Range("D16").Value2 = 0.5200879928886192
$ws2.Range("E16").Value2 = 0.09552581112637204
$ws2.Range("F16").Value2 = -0.1315017241078608
$ws2.Range("A17").Value2 = "QAGEDEP"
$ws2.Range("B17").Value2 = 0.01678028746183512
$ws2.Range("C17").Value2 = -0.1306700641906274
$ws2.Range("D17").Value2 = -0.09147423709580581
$ws2.Range("E17").Value2 = 0.6945701934258486
$ws2.Range("F17").Value2 = 0.5699008906849055
$ws2.Range("A18").Value2 = "QFEMLBR"
$ws2.Range("B18").Value2 = -0.1982536593004132
$ws2.Range("C18").Value2 = 0.1418071566825722
$ws2.Range("D18").Value2 = 0.03868226967556371
$ws2.Range("E18").Value2 = 0.7433016918963765
$ws2.Range("F18").Value2 = -0.02001253281241641
$ws2.Range("A19").Value2 = "QFEMALE"
$ws2.Range("B19").Value2 = -0.0665850047221718
$ws2.Range("C19").Value2 = -0.06807863468096388
$ws2.Range("D19").Value2 = -0.02048003800765474
$ws2.Range("E19").Value2 = 0.8696948779508056
$ws2.Range("F19").Value2 = 0.1203783652956305
$ws2.Range("A20").Value2 = "QSSBEN"
$ws2.Range("B20").Value2 = 0.05321156937084062
$ws2.Range("C20").Value2 = -0.04844388660497572
$ws2.Range("D20").Value2 = -0.1337452429896528
$ws2.Range("E20").Value2 = 0.1584510228281707
$ws2.Range("F20").Value2 = 0.7502919683851625
$ws2.Range("A21").Value2 = "MEDAGE"
$ws2.Range("B21").Value2 = -0.324447414045288
$ws2.Range("C21").Value2 = -0.2578790660438089
$ws2.Range("D21").Value2 = -0.2749783775127918
$ws2.Range("E21").Value2 = -0.001848683687433284
$ws2.Range("F21").Value2 = 0.8112807265660074

# --- Sheet 3: All Refactor Variances ---
$ws3 = $wb.Worksheets.Item("All Refactor Variances")
$ws3.Range("I2").Value2 = 5.029294816679101
$ws3.Range("J2").Value2 = 3.08851212667733
$ws3.Range("K2").Value2 = 2.341179475069111
$ws3.Range("L2").Value2 = 2.102363003151583
$ws3.Range("M2").Value2 = 2.065199329427239
$ws3.Range("N2").Value2 = 5.140423307946821
$ws3.Range("O2").Value2 = 2.750776491517093
$ws3.Range("P2").Value2 = 2.381629080808134
$ws3.Range("Q2").Value2 = 2.101428977954915
$ws3.Range("R2").Value2 = 2.048359019557499
$ws3.Range("I3").Value2 = 0.2394902293656715
$ws3.Range("J3").Value2 = 0.1470720060322538
$ws3.Range("K3").Value2 = 0.1114847369080529
$ws3.Range("L3").Value2 = 0.1001125239595992
$ws3.Range("M3").Value2 = 0.09834282521082088
$ws3.Range("N3").Value2 = 0.257021165397341
$ws3.Range("O3").Value2 = 0.1375388245758547
$ws3.Range("P3").Value2 = 0.1190814540404067
$ws3.Range("Q3").Value2 = 0.1050714488977457
$ws3.Range("R3").Value2 = 0.1024179509778749
$ws3.Range("I4").Value2 = 0.2394902293656715
$ws3.Range("J4").Value2 = 0.3865622353979253
$ws3.Range("K4").Value2 = 0.4980469723059782
$ws3.Range("L4").Value2 = 0.5981594962655774
$ws3.Range("M4").Value2 = 0.6965023214763983
$ws3.Range("N4").Value2 = 0.257021165397341
$ws3.Range("O4").Value2 = 0.3945599899731957
$ws3.Range("P4").Value2 = 0.5136414440136025
$ws3.Range("Q4").Value2 = 0.6187128929113482
$ws3.Range("R4").Value2 = 0.7211308438892232
$ws3.Range("I5").Value2 = 0.3438469937300659
$ws3.Range("J5").Value2 = 0.2111579552534736
$ws3.Range("K5").Value2 = 0.1600636975218333
$ws3.Range("L5").Value2 = 0.1437360951610147
$ws3.Range("M5").Value2 = 0.1411952583336125
$ws3.Range("N5").Value2 = 0.3564140510356862
$ws3.Range("O5").Value2 = 0.1907265869173982
$ws3.Range("P5").Value2 = 0.165131550050159
$ws3.Range("Q5").Value2 = 0.1457037232398374
$ws3.Range("R5").Value2 = 0.142024088756919

# --- Sheet 4: Final Variances ---
$ws4 = $wb.Worksheets.Item("Final Variances")
$ws4.Range("B2").Value2 = 5.140423307946821
$ws4.Range("C2").Value2 = 2.750776491517093
$ws4.Range("D2").Value2 = 2.381629080808134
$ws4.Range("E2").Value2 = 2.101428977954915
$ws4.Range("F2").Value2 = 2.048359019557499
$ws4.Range("B3").Value2 = 0.257021165397341
$ws4.Range("C3").Value2 = 0.1375388245758547
$ws4.Range("D3").Value2 = 0.1190814540404067
$ws4.Range("E3").Value2 = 0.1050714488977457
$ws4.Range("F3").Value2 = 0.1024179509778749
$ws4.Range("B4").Value2 = 0.257021165397341
$ws4.Range("C4").Value2 = 0.3945599899731957
$ws4.Range("D4").Value2 = 0.5136414440136025
$ws4.Range("E4").Value2 = 0.6187128929113482
$ws4.Range("F4").Value2 = 0.7211308438892232
$ws4.Range("B5").Value2 = 0.3564140510356862
$ws4.Range("C5").Value2 = 0.1907265869173982
$ws4.Range("D5").Value2 = 0.165131550050159
$ws4.Range("E5").Value2 = 0.1457037232398374
$ws4.Range("F5").Value2 = 0.142024088756919

# --- Sheet 5: Included and Excluded ---
$ws5 = $wb.Worksheets.Item("Included and Excluded")
$ws5.Range("B2").Value2 = "[['QSERV', 'QFHH', 'QHISPC', 'QEDLESHI', 'QNOHLTH', 'QESL', 'PPUNIT', 'QEXTRCT', 'PERCAP', 'QRICH', 'MDHSEVAL', 'QNOAUTO', 'QPOVTY', 'QRENTER', 'QFAM', 'QAGEDEP', 'QFEMLBR', 'QFEMALE', 'QSSBEN', 'MEDAGE']]"
